$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert the new paragraph right after "Nedan presenteras fynd..."
#    (third paragraph of the body).
# ------------------------------------------------------------------
$introPara = $d.Paragraphs(3)
$introPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(4)
$newPara.Range.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

# ------------------------------------------------------------------
# 2. Remove the old trailing copy of that paragraph together with the
#    two empty paragraphs that preceded it, near the end of the body.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$lastPara.Range.Delete()

# Two empty paragraphs now remain between the last "Kommentar: ..."
# paragraph and the section break. While the last paragraph in the
# document is empty, merge it away by deleting the paragraph mark of
# the paragraph right before it.
while ($true) {
    $n = $d.Paragraphs.Count
    $last = $d.Paragraphs($n)
    $lastLen = $last.Range.Text.TrimEnd([char]13).Length
    if ($lastLen -eq 0 -and $n -gt 1) {
        $prev = $d.Paragraphs($n - 1)
        $mark = $d.Range($prev.Range.End - 1, $prev.Range.End)
        $mark.Delete()
    } else {
        break
    }
}

# ------------------------------------------------------------------
# 3. Bump the date in the first-page header from 2023-11-13 to
#    2023-11-14.
# ------------------------------------------------------------------
$headerRange = $d.Sections(1).Headers(2).Range
$headerRange.Find.Execute("2023-11-13", $true, $false, $false, $false, $false,
                           $true, 1, $false, "2023-11-14", 2)
